# missingVariable-template.docx
#
# Upstream commit: "Fixed #295 Add the version of M2Doc in the template
# custom properties."
#
# The canonical-OOXML diff for *this particular* resource, however, does
# not contain any semantic change: every <w:.../> element in
# word/document.xml and word/styles.xml keeps exactly the same tag name,
# the same set of attribute name/value pairs and the same empty-element
# shape on both sides of the diff - only the *serialized order* of the
# attributes (and of the root element's xmlns declarations) differs,
# which is a non-visible, non-semantic artifact of how the authoring
# tool re-emitted the part (e.g. alphabetised attributes) when the test
# resources were regenerated for the M2Doc-version feature. No text,
# run/paragraph formatting, style definition, section/page setup value,
# or document property actually changed for this file.
#
# So there is nothing for Word's object model to edit here. We simply
# touch the document through a few read-only, side-effect-free OM calls
# (to confirm we are looking at the expected content) and resave it,
# mirroring the harmless resave that produced the upstream attribute
# reordering, without introducing any content change of our own.

$d = $word.ActiveDocument

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
Write-Output ("Tables: " + $d.Tables.Count)
Write-Output ("Sections: " + $d.Sections.Count)
Write-Output $d.Content.Text

$d.Save()
